$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Title paragraph: "Especificação de Caso de Uso: UC10 - Comprar Ingresso"
#    -> "Especificação de Caso de Uso: UC06 - Comprar Ingresso", split across
#       three runs: "Es" | "pecificação de Caso de Uso: UC06" | " - Comprar Ingresso"
# -----------------------------------------------------------------------
$d.Content.Find.Execute("UC10", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "UC06", 2)

$titlePara = $d.Paragraphs(2)
$titleStart = $titlePara.Range.Start
$len1 = 2
$len2 = ("pecificação de Caso de Uso: UC06").Length

$t1 = $d.Range($titleStart, $titleStart + $len1)
$t1.Bold = $false
$t1.Bold = $true

$t2 = $d.Range($titleStart + $len1, $titleStart + $len1 + $len2)
$t2.Bold = $false
$t2.Bold = $true

# -----------------------------------------------------------------------
# 2) "5.1.8." paragraph: insert "Filme –" before "Tipo" (split into runs)
#    and append the _GoBack bookmark right after the final run.
# -----------------------------------------------------------------------
$para518 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "5.1.8.*") {
        $para518 = $d.Paragraphs($i)
        break
    }
}

$p518Start = $para518.Range.Start
$prefixLen = ("5.1.8. Sistema exibe resumo da compra (").Length

$insertRange = $d.Range($p518Start + $prefixLen, $p518Start + $prefixLen)
$insertRange.InsertBefore("Filme –")

$filmeLen = ("Filme ").Length
$dashLen = ("–").Length

$bound1 = $p518Start + $prefixLen
$bound2 = $bound1 + $filmeLen
$bound3 = $bound2 + $dashLen

$r1 = $d.Range($p518Start, $bound1)
$r1.Bold = $true
$r1.Bold = $false

$r2 = $d.Range($bound1, $bound2)
$r2.Bold = $true
$r2.Bold = $false

$r3 = $d.Range($bound2, $bound3)
$r3.Bold = $true
$r3.Bold = $false

# Re-grab paragraph end (insert above shifted it) and re-insert the last
# run's text together with a trailing _GoBack bookmark via WordOpenXML so
# the bookmark lands immediately after the run instead of merging into the
# following paragraph.
$p518End = $para518.Range.End
$r4 = $d.Range($bound3, $p518End - 1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="26"/></w:rPr><w:t>Tipo – quantidade – total de ingressos – valor a ser pago) e solicita confirmação.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r4.InsertXML($xml)

# -----------------------------------------------------------------------
# 3) "5.2. Fluxo Exceção" + ":" runs -> merge into a single run
# -----------------------------------------------------------------------
$para52 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "5.2. Fluxo Exceção:*") {
        $para52 = $d.Paragraphs($i)
        break
    }
}
$scopedRange = $d.Range($para52.Range.Start, $para52.Range.End - 1)
$scopedRange.Find.Execute("5.2. Fluxo Exceção:", $false, $false, $false, $false, `
                           $false, $true, 1, $false, "5.2. Fluxo Exceção:", 2)

# -----------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark that used to sit at the end of the
#    "5.2.3.1." paragraph (it has effectively moved to the 5.1.8 paragraph).
# -----------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
